$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 1805
$ws1.Range("F10").Value = 13
$ws1.Range("F13").Value = 75
$ws1.Range("F15").Value = 701
$ws1.Range("F16").Value = 180
$ws1.Range("F30").Value = 45

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 17

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 1805
$ws4.Range("F12").Value = 13
$ws4.Range("F15").Value = 75
$ws4.Range("F17").Value = 701
$ws4.Range("F18").Value = 180
$ws4.Range("F24").Value = 17
$ws4.Range("F40").Value = 45
